# Auto-generated edit script applying scheduled-runner value updates
# to the Hyperion_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = ""

$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = ""

$ws.Range("H62").Value = 5012.72
$ws.Range("I62").Value = 3707.1428
$ws.Range("K62").Value = 3707.1428
$ws.Range("M62").Value = -3083.1428

$ws.Range("H65").Value = 5012.72
$ws.Range("I65").Value = 3707.1428
$ws.Range("K65").Value = 18535.714
$ws.Range("M65").Value = -15415.714

$ws.Range("H112").Value = 3726.182
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 3789.5813
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 11368.7439
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -13584.7439

$ws.Range("H134").Value = 105882.86
$ws.Range("J134").Value = 105882.86
$ws.Range("L134").Value = 105882.86
$ws.Range("N134").Value = -116022.86

$ws.Range("H137").Value = 61552.633
$ws.Range("I137").Value = 73397.03999999999
$ws.Range("J137").Value = 2330.6
$ws.Range("K137").Value = 220191.12
$ws.Range("L137").Value = 6991.799999999999
$ws.Range("M137").Value = -217641.12
$ws.Range("N137").Value = -12091.8

$ws.Range("H138").Value = 3186.25
$ws.Range("I138").Value = 1221.3334
$ws.Range("J138").Value = 3995.3333
$ws.Range("K138").Value = 3664.0002
$ws.Range("L138").Value = 11985.9999
$ws.Range("M138").Value = 1475.9998
$ws.Range("N138").Value = -22265.9999

$ws.Range("H139").Value = 98638.60000000001
$ws.Range("J139").Value = 103298.25
$ws.Range("L139").Value = 103298.25
$ws.Range("N139").Value = -113578.25

$ws.Range("H140").Value = 124900
$ws.Range("J140").Value = 124900
$ws.Range("L140").Value = 124900
$ws.Range("N140").Value = -135260

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 222.5
$ws.Range("I4").Value = 195
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 195
$ws.Range("L4").Value = 250
$ws.Range("M4").Value = -79
$ws.Range("N4").Value = -482

$ws.Range("H5").Value = 307.57144
$ws.Range("I5").Value = 310.8
$ws.Range("J5").Value = 299.5
$ws.Range("K5").Value = 310.8
$ws.Range("L5").Value = 299.5
$ws.Range("M5").Value = -198.8
$ws.Range("N5").Value = -523.5

$ws.Range("H32").Value = 7910.783
$ws.Range("I32").Value = 4656.456
$ws.Range("J32").Value = 22663.732
$ws.Range("K32").Value = 4656.456
$ws.Range("L32").Value = 22663.732
$ws.Range("M32").Value = -4369.456
$ws.Range("N32").Value = -23237.732

$ws.Range("H74").Value = 52897.914
$ws.Range("I74").Value = 29891.844
$ws.Range("K74").Value = 29891.844
$ws.Range("M74").Value = -29017.844

$ws.Range("H77").Value = 52897.914
$ws.Range("I77").Value = 29891.844
$ws.Range("K77").Value = 149459.22
$ws.Range("M77").Value = -145091.22

$ws.Range("H140").Value = 117598
$ws.Range("J140").Value = 117598
$ws.Range("L140").Value = 117598
$ws.Range("N140").Value = -127958

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 307.57144
$ws.Range("I4").Value = 310.8
$ws.Range("J4").Value = 299.5
$ws.Range("K4").Value = 310.8
$ws.Range("L4").Value = 299.5
$ws.Range("M4").Value = -195.8
$ws.Range("N4").Value = -529.5

$ws.Range("H94").Value = 2328708.8
$ws.Range("I94").Value = 2778724.2
$ws.Range("J94").Value = 14342.857
$ws.Range("K94").Value = 2778724.2
$ws.Range("L94").Value = 14342.857
$ws.Range("M94").Value = -2778273.2
$ws.Range("N94").Value = -15244.857

$ws.Range("H132").Value = 93000
$ws.Range("J132").Value = 93000
$ws.Range("L132").Value = 93000
$ws.Range("N132").Value = -103120

$ws.Range("H134").Value = 2018.2603
$ws.Range("I134").Value = 1052.258
$ws.Range("J134").Value = 7463
$ws.Range("K134").Value = 3156.774
$ws.Range("L134").Value = 22389
$ws.Range("M134").Value = -621.7740000000003
$ws.Range("N134").Value = -27459

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20770.283
$ws.Range("I31").Value = 3300.2
$ws.Range("J31").Value = 43557.348
$ws.Range("K31").Value = 3300.2
$ws.Range("L31").Value = 43557.348
$ws.Range("M31").Value = -3005.2
$ws.Range("N31").Value = -44147.348

$ws.Range("H34").Value = 20770.283
$ws.Range("I34").Value = 3300.2
$ws.Range("J34").Value = 43557.348
$ws.Range("K34").Value = 3300.2
$ws.Range("L34").Value = 43557.348
$ws.Range("M34").Value = -3098.2
$ws.Range("N34").Value = -43961.348

$ws.Range("H93").Value = 35950.8
$ws.Range("I93").Value = 19938.5
$ws.Range("K93").Value = 19938.5
$ws.Range("M93").Value = -18066.5

$ws.Range("H105").Value = 3825.7742
$ws.Range("I105").Value = 3632.5417
$ws.Range("J105").Value = 4488.2856
$ws.Range("K105").Value = 3632.5417
$ws.Range("L105").Value = 4488.2856
$ws.Range("M105").Value = -1885.5417
$ws.Range("N105").Value = -7982.2856

$ws.Range("H135").Value = 148376
$ws.Range("J135").Value = 148376
$ws.Range("L135").Value = 148376
$ws.Range("N135").Value = -158516

$ws.Range("H140").Value = 57747.5
$ws.Range("J140").Value = 57747.5
$ws.Range("L140").Value = 57747.5
$ws.Range("N140").Value = -68107.5

$ws.Range("H141").Value = 43734.223
$ws.Range("J141").Value = 43734.223
$ws.Range("L141").Value = 43734.223
$ws.Range("N141").Value = -54094.223

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2178.64
$ws.Range("I132").Value = 1327.5
$ws.Range("J132").Value = 2964.3076
$ws.Range("K132").Value = 11947.5
$ws.Range("L132").Value = 26678.7684
$ws.Range("M132").Value = -9417.5
$ws.Range("N132").Value = -31738.7684

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1323459
$ws.Range("I97").Value = 2646071.2
$ws.Range("J97").Value = 846.6667
$ws.Range("K97").Value = 2646071.2
$ws.Range("L97").Value = 846.6667
$ws.Range("M97").Value = -2645575.2
$ws.Range("N97").Value = -1838.6667

$ws.Range("H135").Value = 123333
$ws.Range("J135").Value = 123333
$ws.Range("L135").Value = 123333
$ws.Range("N135").Value = -133473

$ws.Range("H139").Value = 60761.87
$ws.Range("J139").Value = 60761.87
$ws.Range("L139").Value = 60761.87
$ws.Range("N139").Value = -71041.87

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1533.579
$ws.Range("I55").Value = 1385.75
$ws.Range("J55").Value = 1787
$ws.Range("K55").Value = 1385.75
$ws.Range("L55").Value = 1787
$ws.Range("M55").Value = -1212.75
$ws.Range("N55").Value = -2133

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 135000
$ws.Range("J141").Value = 135000
$ws.Range("L141").Value = 135000
$ws.Range("N141").Value = -145360
